$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H116").Value = 2400
$ws.Range("I116").Value = 2000
$ws.Range("K116").Value = 2000
$ws.Range("M116").Value = 1442
$ws.Range("H135").Value = 972.5641000000001
$ws.Range("I135").Value = 602.125
$ws.Range("K135").Value = 5419.125
$ws.Range("M135").Value = -2884.125
$ws.Range("H137").Value = 2801.3777
$ws.Range("I137").Value = 2835.2144
$ws.Range("J137").Value = 2745.647
$ws.Range("K137").Value = 8505.643199999999
$ws.Range("L137").Value = 8236.940999999999
$ws.Range("M137").Value = -5955.643199999999
$ws.Range("N137").Value = -13336.941
$ws.Range("H138").Value = 2149.5405
$ws.Range("I138").Value = 1407.5366
$ws.Range("J138").Value = 3071.4243
$ws.Range("K138").Value = 4222.6098
$ws.Range("L138").Value = 9214.2729
$ws.Range("M138").Value = 917.3901999999998
$ws.Range("N138").Value = -19494.2729

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H25").Value = 1900
$ws.Range("I25").Value = 1900
$ws.Range("K25").Value = 1900
$ws.Range("M25").Value = -1498
$ws.Range("H32").Value = 540533.9399999999
$ws.Range("I32").Value = 569943.0600000001
$ws.Range("J32").Value = 40578.5
$ws.Range("K32").Value = 569943.0600000001
$ws.Range("L32").Value = 40578.5
$ws.Range("M32").Value = -569656.0600000001
$ws.Range("N32").Value = -41152.5
$ws.Range("H102").Value = 5202.857
$ws.Range("I102").Value = 5236.6665
$ws.Range("K102").Value = 5236.6665
$ws.Range("M102").Value = -3614.6665
$ws.Range("H122").Value = 2794.1333
$ws.Range("J122").Value = 4325
$ws.Range("L122").Value = 12975
$ws.Range("N122").Value = -17875
$ws.Range("H132").Value = 3159.6155
$ws.Range("I132").Value = 2210.7778
$ws.Range("J132").Value = 5294.5
$ws.Range("K132").Value = 6632.3334
$ws.Range("L132").Value = 15883.5
$ws.Range("M132").Value = -4102.3334
$ws.Range("N132").Value = -20943.5
$ws.Range("H133").Value = 26000
$ws.Range("J133").Value = 26000
$ws.Range("L133").Value = 26000
$ws.Range("N133").Value = -31060

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H46").Value = 62032.5
$ws.Range("J46").Value = 62032.5
$ws.Range("L46").Value = 62032.5
$ws.Range("N46").Value = -62628.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5798.759
$ws.Range("I31").Value = 1330.96
$ws.Range("J31").Value = 9183.454
$ws.Range("K31").Value = 1330.96
$ws.Range("L31").Value = 9183.454
$ws.Range("M31").Value = -1035.96
$ws.Range("N31").Value = -9773.454
$ws.Range("H34").Value = 5798.759
$ws.Range("I34").Value = 1330.96
$ws.Range("J34").Value = 9183.454
$ws.Range("K34").Value = 1330.96
$ws.Range("L34").Value = 9183.454
$ws.Range("M34").Value = -1128.96
$ws.Range("N34").Value = -9587.454
$ws.Range("H43").Value = 191500
$ws.Range("J43").Value = 191500
$ws.Range("L43").Value = 191500
$ws.Range("N43").Value = -191868
$ws.Range("H58").Value = 2000.1111
$ws.Range("I58").Value = 1345.4546
$ws.Range("J58").Value = 3028.8572
$ws.Range("K58").Value = 1345.4546
$ws.Range("L58").Value = 3028.8572
$ws.Range("M58").Value = -1142.4546
$ws.Range("N58").Value = -3434.8572
$ws.Range("H80").Value = 18000
$ws.Range("I80").Value = 18000
$ws.Range("J80").Value = 0
$ws.Range("K80").Value = 18000
$ws.Range("L80").Value = 0
$ws.Range("M80").Value = -16877
$ws.Range("N80").ClearContents()
$ws.Range("H82").Value = 99181
$ws.Range("J82").Value = 99181
$ws.Range("L82").Value = 99181
$ws.Range("N82").Value = -99903
$ws.Range("H83").Value = 18000
$ws.Range("I83").Value = 18000
$ws.Range("J83").Value = 0
$ws.Range("K83").Value = 54000
$ws.Range("L83").Value = 0
$ws.Range("M83").Value = -48384
$ws.Range("N83").ClearContents()
$ws.Range("H85").Value = 99181
$ws.Range("J85").Value = 99181
$ws.Range("L85").Value = 99181
$ws.Range("N85").Value = -101677
$ws.Range("H88").Value = 16666.334
$ws.Range("I88").Value = 0
$ws.Range("J88").Value = 16666.334
$ws.Range("K88").Value = 0
$ws.Range("L88").Value = 16666.334
$ws.Range("M88").ClearContents()
$ws.Range("N88").Value = -17478.334
$ws.Range("H91").Value = 16666.334
$ws.Range("I91").Value = 0
$ws.Range("J91").Value = 16666.334
$ws.Range("K91").Value = 0
$ws.Range("L91").Value = 16666.334
$ws.Range("M91").ClearContents()
$ws.Range("N91").Value = -19474.334
$ws.Range("H101").Value = 191500
$ws.Range("J101").Value = 191500
$ws.Range("L101").Value = 191500
$ws.Range("N101").Value = -197990
$ws.Range("H132").Value = 32408954
$ws.Range("I132").Value = 37038340
$ws.Range("J132").Value = 18520796
$ws.Range("K132").Value = 111115020
$ws.Range("L132").Value = 55562388
$ws.Range("M132").Value = -111112490
$ws.Range("N132").Value = -55567448
$ws.Range("H136").Value = 2000.1111
$ws.Range("I136").Value = 1345.4546
$ws.Range("J136").Value = 3028.8572
$ws.Range("K136").Value = 4036.3638
$ws.Range("L136").Value = 9086.571599999999
$ws.Range("M136").Value = -1486.3638
$ws.Range("N136").Value = -14186.5716

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 4485.467
$ws.Range("I3").Value = 3645.3076
$ws.Range("J3").Value = 9946.5
$ws.Range("K3").Value = 10935.9228
$ws.Range("L3").Value = 29839.5
$ws.Range("M3").Value = -10823.9228
$ws.Range("N3").Value = -30063.5
$ws.Range("H7").Value = 398.16666
$ws.Range("I7").Value = 165.33333
$ws.Range("J7").Value = 631
$ws.Range("K7").Value = 495.99999
$ws.Range("L7").Value = 1893
$ws.Range("M7").Value = -383.99999
$ws.Range("N7").Value = -2117
$ws.Range("H92").Value = 785.4286
$ws.Range("I92").Value = 499.66666
$ws.Range("J92").Value = 999.75
$ws.Range("K92").Value = 1498.99998
$ws.Range("L92").Value = 2999.25
$ws.Range("M92").Value = -250.9999800000001
$ws.Range("N92").Value = -5495.25
$ws.Range("H131").Value = 1548.5625
$ws.Range("J131").Value = 1691.2142
$ws.Range("L131").Value = 5073.642599999999
$ws.Range("N131").Value = -15153.6426

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 2278.5715
$ws.Range("I113").Value = 2324.8125
$ws.Range("J113").Value = 2130.6
$ws.Range("K113").Value = 2324.8125
$ws.Range("L113").Value = 2130.6
$ws.Range("M113").Value = -154.8125
$ws.Range("N113").Value = -6470.6
$ws.Range("H126").Value = 2750
$ws.Range("I126").Value = 2500
$ws.Range("J126").Value = 3000
$ws.Range("K126").Value = 7500
$ws.Range("L126").Value = 9000
$ws.Range("M126").Value = -5030
$ws.Range("N126").Value = -13940
$ws.Range("H132").Value = 3539.3845
$ws.Range("I132").Value = 3277.9
$ws.Range("J132").Value = 4411
$ws.Range("K132").Value = 9833.700000000001
$ws.Range("L132").Value = 13233
$ws.Range("M132").Value = -7303.700000000001
$ws.Range("N132").Value = -18293

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H9").Value = 70007
$ws.Range("J9").Value = 70007
$ws.Range("L9").Value = 70007
$ws.Range("N9").Value = -70287
$ws.Range("H14").Value = 852000.4399999999
$ws.Range("I14").Value = 38000
$ws.Range("J14").Value = 1014800.5
$ws.Range("K14").Value = 38000
$ws.Range("L14").Value = 1014800.5
$ws.Range("M14").Value = -37832
$ws.Range("N14").Value = -1015136.5
$ws.Range("H15").Value = 44003.5
$ws.Range("J15").Value = 44003.5
$ws.Range("L15").Value = 44003.5
$ws.Range("N15").Value = -44579.5
$ws.Range("H20").Value = 23813.875
$ws.Range("I20").Value = 0
$ws.Range("J20").Value = 23813.875
$ws.Range("K20").Value = 0
$ws.Range("L20").Value = 23813.875
$ws.Range("M20").ClearContents()
$ws.Range("N20").Value = -24293.875
$ws.Range("H21").Value = 69666.664
$ws.Range("J21").Value = 80000
$ws.Range("L21").Value = 80000
$ws.Range("N21").Value = -80470
$ws.Range("H35").Value = 69666.664
$ws.Range("J35").Value = 80000
$ws.Range("L35").Value = 80000
$ws.Range("N35").Value = -80580
$ws.Range("H54").Value = 19799
$ws.Range("J54").Value = 19799
$ws.Range("L54").Value = 19799
$ws.Range("N54").Value = -20839
$ws.Range("H81").Value = 3168.389
$ws.Range("I81").Value = 3119.5293
$ws.Range("J81").Value = 3999
$ws.Range("K81").Value = 6239.0586
$ws.Range("L81").Value = 7998
$ws.Range("M81").Value = -5178.0586
$ws.Range("N81").Value = -10120
$ws.Range("H84").Value = 3168.389
$ws.Range("I84").Value = 3119.5293
$ws.Range("J84").Value = 3999
$ws.Range("K84").Value = 31195.293
$ws.Range("L84").Value = 39990
$ws.Range("M84").Value = -25891.293
$ws.Range("N84").Value = -50598
$ws.Range("H96").Value = 3679.818
$ws.Range("I96").Value = 1922.25
$ws.Range("J96").Value = 8366.666999999999
$ws.Range("K96").Value = 1922.25
$ws.Range("L96").Value = 8366.666999999999
$ws.Range("M96").Value = -549.25
$ws.Range("N96").Value = -11112.667
$ws.Range("H132").Value = 4275791.5
$ws.Range("I132").Value = 2277.2856
$ws.Range("J132").Value = 15153828
$ws.Range("K132").Value = 6831.8568
$ws.Range("L132").Value = 45461484
$ws.Range("M132").Value = -4301.8568
$ws.Range("N132").Value = -45466544
$ws.Range("H136").Value = 2544.225
$ws.Range("I136").Value = 2068.5833
$ws.Range("J136").Value = 6825
$ws.Range("K136").Value = 6205.749899999999
$ws.Range("L136").Value = 20475
$ws.Range("M136").Value = -3655.749899999999
$ws.Range("N136").Value = -25575
